# Applies the "Added a few more slots" edit:
#  1. Inserts a new "Meta description: ..." paragraph right after the H1 title.
#  2. Removes the duplicated bold title paragraph near the end of the document.
#  3. Rewrites the trailing italic paragraph with the new DALLE image-prompt text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: insert the "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$metaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Crystal Quest Frostlands features action-packed adventure, cascading reels and free spins. Play for free to enjoy the game''s stunning design and big winning opportunities.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaRange.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# Part 2: find the duplicated bold title paragraph near the end (skip the
# real H1 at Paragraphs(1), which has identical text) and the italic
# paragraph that immediately follows it, then replace both paragraphs with a
# single paragraph carrying the new italic text.
# ---------------------------------------------------------------------------
$oldHeadingText = 'Play Crystal Quest Frostlands Free - Exciting Slot Game'
$oldItalicText = 'Crystal Quest Frostlands features action-packed adventure, cascading reels and free spins. Play for free to enjoy the game''s stunning design and big winning opportunities.'

$count = $d.Paragraphs.Count
$boldDupe = $null
for ($i = $count; $i -ge 2; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd("`r") -eq $oldHeadingText) {
        $boldDupe = $p
        break
    }
}

if ($boldDupe -ne $null) {
    $italicPara = $boldDupe.Next()
    if ($italicPara -ne $null -and $italicPara.Range.Text.TrimEnd("`r") -eq $oldItalicText) {
        $replaceRange = $d.Range($boldDupe.Range.Start, $italicPara.Range.End)
        $bodyXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Create a feature image fitting Crystal Quest Frostlands: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses DALLE, can you create a feature image for Crystal Quest Frostlands? The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be standing in a snowy landscape, holding a crystal and fighting against a demonic yeti. The image should also include the game''s logo, "Crystal Quest Frostlands," in bold letters and bright colors. The background should be icy with snowflakes falling, and the overall tone should be adventurous and exciting. Let''s capture the thrill of the game with a visually stunning, attention-grabbing feature image!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $replaceRange.InsertXML($bodyXml)
    }
}

Write-Output "edit complete"
